$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Names (column A), entered in row order ---
$ws.Range("A4").Value = "Gustavo"
$ws.Range("A5").Value = "Guilherme"
$ws.Range("A6").Value = "Thiago"
$ws.Range("A7").Value = "Andreia"
$ws.Range("A8").Value = "Pguedes"

# --- Emails (column B), entered in the same order the author typed them ---
$ws.Range("B5").Value = "guipsguedes@gmail.com"
$ws.Range("B7").Value = "art.andreiaguedes@gmail.com"
$ws.Range("B8").Value = "paulogsguedes@gmail.com"
$ws.Range("B4").Value = "gpsguedes@gmail.com"
$ws.Range("B6").Value = "tpsguedes@gmail.com"

# --- Turn each new email into a live mailto hyperlink (same as B2/B3) ---
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:gpsguedes@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:guipsguedes@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:art.andreiaguedes@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:paulogsguedes@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:tpsguedes@gmail.com")

# Hyperlinks.Add() re-styles the cell with a fresh xf; put it back to the
# plain hyperlink style already used by the existing B2/B3 cells.
$ws.Range("B4:B8").Style = $ws.Range("B2").Style

# Final selection left on the header row, matching the saved workbook.
$ws.Range("A1:B1").Select() | Out-Null
